$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Docente\(s\) Respons") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph 'Docente(s) Responsável(eis)'"
}

# Insert a new (initially empty) paragraph right after it
$target.Range.InsertParagraphAfter()
$newPara = $target.Next()

# Build the new paragraph's content (a ListBullet-styled paragraph with three
# runs, the first two each followed by a manual line break) via a small
# WordprocessingML package fragment inserted through Range.InsertXML so that
# each name ends up as its own separate run, matching how Word itself
# represents text typed with Shift+Enter between runs of differing history.
$xml = @'
<?xml version="1.0"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml">
    <pkg:xmlData>
      <Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">
        <Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>
      </Relationships>
    </pkg:xmlData>
  </pkg:part>
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListBullet"/>
            </w:pPr>
            <w:r>
              <w:t>7459752 - Maria Ismenia Sodero Toledo Faria</w:t>
              <w:br/>
            </w:r>
            <w:r>
              <w:t>2166002 - Sandra Giacomin Schneider</w:t>
              <w:br/>
            </w:r>
            <w:r>
              <w:t>1922320 - Sebastiao Ribeiro</w:t>
            </w:r>
          </w:p>
          <w:sectPr/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($xml)
